$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.511.54"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "2.046.14"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'246.89"
$ws.Range("D5").Style = $ws.Range("D4").Style
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'0.659"
$ws.Range("D6").Style = $ws.Range("D4").Style
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'54.32"
$ws.Range("D8").Style = $ws.Range("D4").Style
$ws.Range("E8").Value = "  -6.49%  "
$ws.Range("D9").Value = "'62.90"
$ws.Range("D9").Style = $ws.Range("D4").Style
$ws.Range("E9").Value = "  +5.67%  "
$ws.Range("D10").Value = "'0.363"
$ws.Range("D10").Style = $ws.Range("D4").Style
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").Value = "'0.0746"
$ws.Range("D11").Style = $ws.Range("D4").Style
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("E12").Value = "  -4.04%  "
$ws.Range("D13").Value = "'0.937"
$ws.Range("D13").Style = $ws.Range("D4").Style
$ws.Range("E13").Value = "  +6.25%  "
$ws.Range("D14").Value = "'14.59"
$ws.Range("D14").Style = $ws.Range("D4").Style
$ws.Range("E14").Value = "  -4.82%  "
$ws.Range("D15").Value = "2.342.59"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'5.42"
$ws.Range("D16").Style = $ws.Range("D4").Style
$ws.Range("E16").Value = "  -4.84%  "
$ws.Range("D17").Value = "2.045.51"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "36.357.26"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'17.01"
$ws.Range("D19").Style = $ws.Range("D4").Style
$ws.Range("E19").Value = "  -5.70%  "
$ws.Range("D20").Value = "'71.61"
$ws.Range("D20").Style = $ws.Range("D4").Style
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "0.0₃0855"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").Value = "'236.68"
$ws.Range("D22").Style = $ws.Range("D4").Style
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'5.18"
$ws.Range("D23").Style = $ws.Range("D4").Style
$ws.Range("E23").Value = "  -4.80%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("D26").Style = $ws.Range("D4").Style
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").Value = "'164.39"
$ws.Range("D27").Style = $ws.Range("D4").Style
$ws.Range("E27").Value = "  -3.12%  "
$ws.Range("E28").Value = "  -13.03%  "
$ws.Range("D29").Value = "'19.87"
$ws.Range("D29").Style = $ws.Range("D4").Style
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D31").Value = "'5.03"
$ws.Range("D31").Style = $ws.Range("D4").Style
$ws.Range("E31").Value = "  -9.35%  "
$ws.Range("D32").Value = "'1.18"
$ws.Range("D32").Style = $ws.Range("D4").Style
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("D34").Value = "'4.38"
$ws.Range("D34").Style = $ws.Range("D4").Style
$ws.Range("E34").Value = "  -7.66%  "
$ws.Range("D35").Value = "'0.0880"
$ws.Range("D35").Style = $ws.Range("D4").Style
$ws.Range("E35").Value = "  +7.34%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").Value = "'2.21"
$ws.Range("D38").Style = $ws.Range("D4").Style
$ws.Range("E38").Value = "  -6.65%  "
$ws.Range("D39").Value = "'5.02"
$ws.Range("D39").Style = $ws.Range("D4").Style
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").Value = "'1.23"
$ws.Range("D40").Style = $ws.Range("D4").Style
$ws.Range("E40").Value = "  -7.49%  "
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("E42").Value = "  -5.24%  "
$ws.Range("E43").Value = "  -4.85%  "
$ws.Range("D44").Value = "'93.52"
$ws.Range("D44").Style = $ws.Range("D4").Style
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("D45").Value = "'0.0897"
$ws.Range("D45").Style = $ws.Range("D4").Style
$ws.Range("E45").Value = "  -6.62%  "
$ws.Range("D46").Value = "1.390.45"
$ws.Range("E46").Value = "  +6.11%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.43"
$ws.Range("D47").Style = $ws.Range("D4").Style
$ws.Range("E47").Value = "  +9.47%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'15.64"
$ws.Range("D48").Style = $ws.Range("D4").Style
$ws.Range("E48").Value = "  -8.19%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").Value = "  -5.57%  "
$ws.Range("D51").Value = "2.226.91"
$ws.Range("E51").Value = "  -0.60%  "
